# Remove the two test rows ("AaTest" / "BBTest") that were temporarily added
# near the top of the data table. Everything below shifts up by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("3:4").Delete()

# The worksheet keeps a sortState (ref="A5:D44" before the delete) describing
# the range that was last sorted. Deleting rows doesn't retarget it on its
# own, so it would be left stale (still "A5:D44") after the shift. Re-apply
# the existing order (via a throwaway sequential key in column A, restored
# right after) so the sort machinery re-anchors sortState to the new
# "A3:D42" range without actually reordering any rows.
$firstRow = 3
$lastRow = 42
$original = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $original[$r] = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $r
}

$sortObj = $ws.Sort
$sortFields = $sortObj.SortFields
$sortFields.Clear()
$sortFields.Add($ws.Range("A" + $firstRow + ":A" + $lastRow)) | Out-Null
$sortObj.SetRange($ws.Range("A" + $firstRow + ":D" + $lastRow)) | Out-Null
$sortObj.Header = -4135
$sortObj.Apply()

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $original[$r]
}

# Match the author's final selection state.
$ws.Range("A3:XFD4").Select() | Out-Null
